$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.36869142360349
$ws.Range("C2").Value = 10.75117040981646
$ws.Range("E2").Value = 16.59640256886743
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 27.77976244798875
$ws.Range("H2").Value = 12.93568528974765
$ws.Range("N2").Value = 15.73962196915567

$ws.Range("B3").Value = 13.62612902144914
$ws.Range("C3").Value = 10.1408211776373
$ws.Range("E3").Value = 15.64447618581514
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 27.10757498782016
$ws.Range("H3").Value = 12.93453169931673
$ws.Range("N3").Value = 15.81610850659851

$ws.Range("B4").Value = 13.15196288439352
$ws.Range("C4").Value = 9.745113458587241
$ws.Range("E4").Value = 15.03507972112006
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 26.70278110408819
$ws.Range("H4").Value = 12.93829847740028
$ws.Range("N4").Value = 15.8650952768412

$ws.Range("B5").Value = 12.95439730688211
$ws.Range("C5").Value = 9.57864225974874
$ws.Range("E5").Value = 14.78074912488671
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 26.54011844935768
$ws.Range("H5").Value = 12.94095267276592
$ws.Range("N5").Value = 15.88556893154309

$ws.Range("B6").Value = 12.92133761591868
$ws.Range("C6").Value = 9.550686507844027
$ws.Range("E6").Value = 14.73816440761065
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 26.51325573325215
$ws.Range("H6").Value = 12.94146080071145
$ws.Range("N6").Value = 15.88899950483655

$ws.Range("B7").Value = 13.14931565996324
$ws.Range("C7").Value = 9.742889415983628
$ws.Range("E7").Value = 15.03167361710603
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 26.70057769478683
$ws.Range("H7").Value = 12.9383297500927
$ws.Range("N7").Value = 15.86536931917056

$ws.Range("B8").Value = 14.11657396044475
$ws.Range("C8").Value = 10.54509702556265
$ws.Range("E8").Value = 16.273492889737
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 27.54652562113571
$ws.Range("H8").Value = 12.93435533977691
$ws.Range("N8").Value = 15.7655758720765

$ws.Range("B9").Value = 15.85991241313657
$ws.Range("C9").Value = 11.95036097894926
$ws.Range("E9").Value = 18.62295077285343
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 29.25529270599567
$ws.Range("H9").Value = 12.96229883117237
$ws.Range("N9").Value = 15.58583406282202

$ws.Range("B10").Value = 17.03772831059452
$ws.Range("C10").Value = 12.87899587172154
$ws.Range("E10").Value = 20.29158400496897
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 30.52407685834825
$ws.Range("H10").Value = 13.00486537700675
$ws.Range("N10").Value = 15.46335501083807

$ws.Range("B11").Value = 17.54973329332596
$ws.Range("C11").Value = 13.27875560321498
$ws.Range("E11").Value = 21.00860635402501
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 31.10078580379488
$ws.Range("H11").Value = 13.02904705215661
$ws.Range("N11").Value = 15.40968394428385

$ws.Range("B12").Value = 17.7401037811087
$ws.Range("C12").Value = 13.42686602129808
$ws.Range("E12").Value = 21.27411787350042
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 31.3188406690878
$ws.Range("H12").Value = 13.03889769083037
$ws.Range("N12").Value = 15.38965179552114

$ws.Range("B13").Value = 17.69926170962584
$ws.Range("C13").Value = 13.39511331015315
$ws.Range("E13").Value = 21.21720189325843
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 31.27189747149325
$ws.Range("H13").Value = 13.03674533269655
$ws.Range("N13").Value = 15.39395312874393

$ws.Range("B14").Value = 17.56546610058793
$ws.Range("C14").Value = 13.29100630283654
$ws.Range("E14").Value = 21.0305703697433
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 31.1187334266611
$ws.Range("H14").Value = 13.02984356966096
$ws.Range("N14").Value = 15.40803004936855

$ws.Range("B15").Value = 17.48305223702924
$ws.Range("C15").Value = 13.22681174187975
$ws.Range("E15").Value = 20.9154716840706
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 31.02486517286586
$ws.Range("H15").Value = 13.02570636661894
$ws.Range("N15").Value = 15.41669052137577

$ws.Range("B16").Value = 17.00378043585164
$ws.Range("C16").Value = 12.85241340943024
$ws.Range("E16").Value = 20.24388198160525
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 30.48635512542668
$ws.Range("H16").Value = 13.00338208572572
$ws.Range("N16").Value = 15.46690351013838

$ws.Range("B17").Value = 16.70359695155069
$ws.Range("C17").Value = 12.61691554278153
$ws.Range("E17").Value = 19.82114056086633
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 30.15568307954285
$ws.Range("H17").Value = 12.99092172955709
$ws.Range("N17").Value = 15.49822983610884

$ws.Range("B18").Value = 16.52870449928729
$ws.Range("C18").Value = 12.47932873859463
$ws.Range("E18").Value = 19.5740309815304
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 29.96546012099068
$ws.Range("H18").Value = 12.98420837384302
$ws.Range("N18").Value = 15.5164405510988

$ws.Range("B19").Value = 16.46910824647594
$ws.Range("C19").Value = 12.43237791380189
$ws.Range("E19").Value = 19.489682676596
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 29.90105724067165
$ws.Range("H19").Value = 12.98201319732342
$ws.Range("N19").Value = 15.52263953924338

$ws.Range("B20").Value = 16.73578403237565
$ws.Range("C20").Value = 12.64220569732863
$ws.Range("E20").Value = 19.86655167175908
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 30.19088867332972
$ws.Range("H20").Value = 12.99220121015962
$ws.Range("N20").Value = 15.49487517351022

$ws.Range("B21").Value = 17.60486115982197
$ws.Range("C21").Value = 13.3216738440904
$ws.Range("E21").Value = 21.08555136785519
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 31.1637324663969
$ws.Range("H21").Value = 13.03185196105669
$ws.Range("N21").Value = 15.40388741134859

$ws.Range("B22").Value = 18.15233567922065
$ws.Range("C22").Value = 13.74668505943578
$ws.Range("E22").Value = 21.84722852666485
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 31.7974944742342
$ws.Range("H22").Value = 13.06180821852032
$ws.Range("N22").Value = 15.34612216657074

$ws.Range("B23").Value = 17.86204203341495
$ws.Range("C23").Value = 13.52159427752354
$ws.Range("E23").Value = 21.44389826995025
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 31.4595142721866
$ws.Range("H23").Value = 13.04545017083939
$ws.Range("N23").Value = 15.3767976887836

$ws.Range("B24").Value = 16.72123945021344
$ws.Range("C24").Value = 12.63077885764615
$ws.Range("E24").Value = 19.84603398817861
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 30.1749725628801
$ws.Range("H24").Value = 12.99162135456144
$ws.Range("N24").Value = 15.49639119088059

$ws.Range("B25").Value = 15.40581300347445
$ws.Range("C25").Value = 11.58835522365276
$ws.Range("E25").Value = 17.97092871007376
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 28.78948437314263
$ws.Range("H25").Value = 12.95088359688785
$ws.Range("N25").Value = 15.63276645064548
